$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1874.2632
$ws.Range("I33").Value = 1927.2
$ws.Range("K33").Value = 1927.2
$ws.Range("M33").Value = -1698.2
$ws.Range("H41").Value = 2761.2104
$ws.Range("I41").Value = 2813.4614
$ws.Range("J41").Value = 2648
$ws.Range("K41").Value = 2813.4614
$ws.Range("L41").Value = 2648
$ws.Range("M41").Value = -2373.4614
$ws.Range("N41").Value = -3528
$ws.Range("H53").Value = 771.1111
$ws.Range("I53").Value = 300
$ws.Range("J53").Value = 905.7143
$ws.Range("K53").Value = 300
$ws.Range("L53").Value = 905.7143
$ws.Range("M53").Value = 337
$ws.Range("N53").Value = -2179.7143
$ws.Range("H55").Value = 1881.9166
$ws.Range("I55").Value = 252
$ws.Range("K55").Value = 252
$ws.Range("M55").Value = -38
$ws.Range("H62").Value = 13039.4
$ws.Range("I62").Value = 4461.5
$ws.Range("J62").Value = 18758
$ws.Range("K62").Value = 4461.5
$ws.Range("L62").Value = 18758
$ws.Range("M62").Value = -3837.5
$ws.Range("N62").Value = -20006
$ws.Range("H65").Value = 13039.4
$ws.Range("I65").Value = 4461.5
$ws.Range("J65").Value = 18758
$ws.Range("K65").Value = 22307.5
$ws.Range("L65").Value = 93790
$ws.Range("M65").Value = -19187.5
$ws.Range("N65").Value = -100030
$ws.Range("H76").Value = 3950
$ws.Range("I76").Value = 4142.857
$ws.Range("K76").Value = 4142.857
$ws.Range("M76").Value = -3827.857
$ws.Range("H79").Value = 3950
$ws.Range("I79").Value = 4142.857
$ws.Range("K79").Value = 4142.857
$ws.Range("M79").Value = -3050.857
$ws.Range("H135").Value = 3102.9412
$ws.Range("I135").Value = 3102.9412
$ws.Range("K135").Value = 27926.4708
$ws.Range("M135").Value = -25391.4708
$ws.Range("H137").Value = 3249.4614
$ws.Range("I137").Value = 2273.5518
$ws.Range("J137").Value = 6079.6
$ws.Range("K137").Value = 6820.655400000001
$ws.Range("L137").Value = 18238.8
$ws.Range("M137").Value = -4270.655400000001
$ws.Range("N137").Value = -23338.8
$ws.Range("H138").Value = 2615.49
$ws.Range("I138").Value = 1394.1072
$ws.Range("J138").Value = 3090.4722
$ws.Range("K138").Value = 4182.321599999999
$ws.Range("L138").Value = 9271.4166
$ws.Range("M138").Value = 957.6784000000007
$ws.Range("N138").Value = -19551.4166

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5683726.5
$ws.Range("I32").Value = 5748596.5
$ws.Range("K32").Value = 5748596.5
$ws.Range("M32").Value = -5748309.5
$ws.Range("H61").Value = 10898202
$ws.Range("I61").Value = 16132669
$ws.Range("J61").Value = 80304.266
$ws.Range("K61").Value = 16132669
$ws.Range("L61").Value = 80304.266
$ws.Range("M61").Value = -16132457
$ws.Range("N61").Value = -80728.266
$ws.Range("H74").Value = 8628738
$ws.Range("I74").Value = 10419029
$ws.Range("J74").Value = 35337
$ws.Range("K74").Value = 10419029
$ws.Range("L74").Value = 35337
$ws.Range("M74").Value = -10418155
$ws.Range("N74").Value = -37085
$ws.Range("H76").Value = 22249
$ws.Range("J76").Value = 22249
$ws.Range("L76").Value = 22249
$ws.Range("N76").Value = -22925
$ws.Range("H77").Value = 8628738
$ws.Range("I77").Value = 10419029
$ws.Range("J77").Value = 35337
$ws.Range("K77").Value = 52095145
$ws.Range("L77").Value = 176685
$ws.Range("M77").Value = -52090777
$ws.Range("N77").Value = -185421
$ws.Range("H79").Value = 22249
$ws.Range("J79").Value = 22249
$ws.Range("L79").Value = 22249
$ws.Range("N79").Value = -24589
$ws.Range("H132").Value = 3937.4075
$ws.Range("I132").Value = 1721.4
$ws.Range("J132").Value = 6707.4165
$ws.Range("K132").Value = 5164.200000000001
$ws.Range("L132").Value = 20122.2495
$ws.Range("M132").Value = -2634.200000000001
$ws.Range("N132").Value = -25182.2495
$ws.Range("H136").Value = 10898202
$ws.Range("I136").Value = 16132669
$ws.Range("J136").Value = 80304.266
$ws.Range("K136").Value = 48398007
$ws.Range("L136").Value = 240912.798
$ws.Range("M136").Value = -48395457
$ws.Range("N136").Value = -246012.798

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2623.8975
$ws.Range("I20").Value = 3003.4062
$ws.Range("J20").Value = 889
$ws.Range("K20").Value = 3003.4062
$ws.Range("L20").Value = 889
$ws.Range("M20").Value = -2756.4062
$ws.Range("N20").Value = -1383
$ws.Range("H86").Value = 1699.6207
$ws.Range("I86").Value = 1224.5834
$ws.Range("K86").Value = 1224.5834
$ws.Range("M86").Value = -101.5834
$ws.Range("H89").Value = 1699.6207
$ws.Range("I89").Value = 1224.5834
$ws.Range("K89").Value = 6122.916999999999
$ws.Range("M89").Value = -506.9169999999995
$ws.Range("H105").Value = 2910.1667
$ws.Range("I105").Value = 2861.25
$ws.Range("K105").Value = 2861.25
$ws.Range("M105").Value = -1114.25
$ws.Range("H133").Value = 40499.75
$ws.Range("I133").Value = 20000
$ws.Range("K133").Value = 20000
$ws.Range("M133").Value = -14940

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1188.8
$ws.Range("I22").Value = 861.25
$ws.Range("K22").Value = 861.25
$ws.Range("M22").Value = -511.25
$ws.Range("H31").Value = 493128.97
$ws.Range("I31").Value = 9130.916999999999
$ws.Range("J31").Value = 977127.0600000001
$ws.Range("K31").Value = 9130.916999999999
$ws.Range("L31").Value = 977127.0600000001
$ws.Range("M31").Value = -8835.916999999999
$ws.Range("N31").Value = -977717.0600000001
$ws.Range("H34").Value = 493128.97
$ws.Range("I34").Value = 9130.916999999999
$ws.Range("J34").Value = 977127.0600000001
$ws.Range("K34").Value = 9130.916999999999
$ws.Range("L34").Value = 977127.0600000001
$ws.Range("M34").Value = -8928.916999999999
$ws.Range("N34").Value = -977531.0600000001
$ws.Range("H107").Value = 947
$ws.Range("I107").Value = 447.83334
$ws.Range("J107").Value = 1612.5555
$ws.Range("K107").Value = 447.83334
$ws.Range("L107").Value = 1612.5555
$ws.Range("M107").Value = 1472.16666
$ws.Range("N107").Value = -5452.5555
$ws.Range("H121").Value = 55000
$ws.Range("J121").Value = 55000
$ws.Range("L121").Value = 55000
$ws.Range("N121").Value = -57620
$ws.Range("H132").Value = 2003.44
$ws.Range("I132").Value = 1677.1951
$ws.Range("K132").Value = 5031.5853
$ws.Range("M132").Value = -2501.5853

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 286.92856
$ws.Range("J7").Value = 269.75
$ws.Range("L7").Value = 809.25
$ws.Range("N7").Value = -1033.25
$ws.Range("H54").Value = 17499.75
$ws.Range("J54").Value = 18333
$ws.Range("L54").Value = 54999
$ws.Range("N54").Value = -56117

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4252.5
$ws.Range("I80").Value = 4505
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 4505
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -3507
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 4252.5
$ws.Range("I83").Value = 4505
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 22525
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -17533
$ws.Range("N83").Value = -29984
$ws.Range("H99").Value = 25947.25
$ws.Range("J99").Value = 41990
$ws.Range("L99").Value = 41990
$ws.Range("N99").Value = -46482
$ws.Range("H128").Value = 84665.55499999999
$ws.Range("J128").Value = 84665.55499999999
$ws.Range("L128").Value = 84665.55499999999
$ws.Range("N128").Value = -94625.55499999999
$ws.Range("H132").Value = 33338902
$ws.Range("I132").Value = 52633070
$ws.Range("K132").Value = 157899210
$ws.Range("M132").Value = -157896680

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 933.8
$ws.Range("I16").Value = 714.7857
$ws.Range("K16").Value = 714.7857
$ws.Range("M16").Value = -544.7857
$ws.Range("H22").Value = 3456.2856
$ws.Range("I22").Value = 3449
$ws.Range("J22").Value = 3500
$ws.Range("K22").Value = 3449
$ws.Range("L22").Value = 3500
$ws.Range("M22").Value = -3154
$ws.Range("N22").Value = -4090
$ws.Range("H27").Value = 3456.2856
$ws.Range("I27").Value = 3449
$ws.Range("J27").Value = 3500
$ws.Range("K27").Value = 3449
$ws.Range("L27").Value = 3500
$ws.Range("M27").Value = -3342
$ws.Range("N27").Value = -3714
$ws.Range("H46").Value = 6430.4546
$ws.Range("I46").Value = 1705
$ws.Range("J46").Value = 14700
$ws.Range("K46").Value = 1705
$ws.Range("L46").Value = 14700
$ws.Range("M46").Value = -1517
$ws.Range("N46").Value = -15076
$ws.Range("H55").Value = 52632056
$ws.Range("J55").Value = 509.25
$ws.Range("L55").Value = 509.25
$ws.Range("N55").Value = -855.25
$ws.Range("H68").Value = 2833
$ws.Range("I68").Value = 2749.5
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 2749.5
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -2000.5
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 2833
$ws.Range("I71").Value = 2749.5
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 13747.5
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -10003.5
$ws.Range("N71").Value = -22488
$ws.Range("H136").Value = 43879.9
$ws.Range("I136").Value = 6760.3887
$ws.Range("K136").Value = 20281.1661
$ws.Range("M136").Value = -17731.1661

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2417.1667
$ws.Range("I100").Value = 2417.1667
$ws.Range("K100").Value = 4834.3334
$ws.Range("M100").Value = -4293.3334
$ws.Range("H126").Value = 3785.4722
$ws.Range("I126").Value = 3926.862
$ws.Range("K126").Value = 11780.586
$ws.Range("M126").Value = -9310.585999999999
$ws.Range("H132").Value = 2403.56
$ws.Range("I132").Value = 2403.7083
$ws.Range("K132").Value = 7211.124899999999
$ws.Range("M132").Value = -4681.124899999999
$ws.Range("H136").Value = 7049.625
$ws.Range("I136").Value = 1233.3334
$ws.Range("K136").Value = 3700.0002
$ws.Range("M136").Value = -1150.0002
